# Add implementation of MSM (Method Size Measure) to the methodNumberOfLines sheet.
#
# Adds constructor rows (Product(), Product(args), StockAppTest(), StockComponentTests(),
# OrderManageService(args), KafkaContainerDevMode(), StockApp()) that were previously
# missing from the method-line-count table, which pushes the table from 19 to 27 data
# rows and reorders the existing rows to keep rows grouped by declaring class.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

# Final (Class Name, Method Signature, Number of Lines) content for rows 2..27, in order.
$data = @(
    @(2, 'pl.piomin.stock.domain.Product', 'Product()', '2'),
    @(3, 'pl.piomin.stock.domain.Product', 'Product(java.lang.Long, java.lang.String, int, int)', '6'),
    @(4, 'pl.piomin.stock.domain.Product', 'getId()', '3'),
    @(5, 'pl.piomin.stock.domain.Product', 'setId(java.lang.Long)', '3'),
    @(6, 'pl.piomin.stock.domain.Product', 'getName()', '3'),
    @(7, 'pl.piomin.stock.domain.Product', 'setName(java.lang.String)', '3'),
    @(8, 'pl.piomin.stock.domain.Product', 'getAvailableItems()', '3'),
    @(9, 'pl.piomin.stock.domain.Product', 'setAvailableItems(int)', '3'),
    @(10, 'pl.piomin.stock.domain.Product', 'getReservedItems()', '3'),
    @(11, 'pl.piomin.stock.domain.Product', 'setReservedItems(int)', '3'),
    @(12, 'pl.piomin.stock.domain.Product', 'toString()', '3'),
    @(13, 'pl.piomin.stock.StockAppTest', 'StockAppTest()', '1'),
    @(14, 'pl.piomin.stock.StockAppTest', 'main(java.lang.String[])', '3'),
    @(15, 'pl.piomin.stock.StockComponentTests', 'StockComponentTests()', '1'),
    @(16, 'pl.piomin.stock.StockComponentTests', 'eventAccept()', '11'),
    @(17, 'pl.piomin.stock.StockComponentTests', 'eventReject()', '10'),
    @(18, 'pl.piomin.stock.StockComponentTests', 'eventConfirm()', '10'),
    @(19, 'pl.piomin.stock.service.OrderManageService', 'OrderManageService(pl.piomin.stock.repository.ProductRepository, org.springframework.kafka.core.KafkaTemplate)', '4'),
    @(20, 'pl.piomin.stock.service.OrderManageService', 'reserve(pl.piomin.base.domain.Order)', '17'),
    @(21, 'pl.piomin.stock.service.OrderManageService', 'confirm(pl.piomin.base.domain.Order)', '13'),
    @(22, 'pl.piomin.stock.KafkaContainerDevMode', 'KafkaContainerDevMode()', '1'),
    @(23, 'pl.piomin.stock.KafkaContainerDevMode', 'kafka()', '3'),
    @(24, 'pl.piomin.stock.StockApp', 'StockApp()', '1'),
    @(25, 'pl.piomin.stock.StockApp', 'main(java.lang.String[])', '3'),
    @(26, 'pl.piomin.stock.StockApp', 'onEvent(pl.piomin.base.domain.Order)', '9'),
    @(27, 'pl.piomin.stock.StockApp', 'generateData()', '8')
)

foreach ($row in $data) {
    $r = $row[0]
    $className = $row[1]
    $methodSig = $row[2]
    $numLines = $row[3]

    $ws.Cells.Item($r, 1).Value = $className
    $ws.Cells.Item($r, 2).Value = $methodSig
    # Prefix with an apostrophe so the numeric-looking line count is stored as text
    # (matching the existing "Number of Lines" column, which is text everywhere else).
    $ws.Cells.Item($r, 3).Value = "'" + $numLines
}
